# "Generate Report for Archive"
# Updates the localization-status report:
#   1. Status text "Ready for handoff" -> "In Translation" for every
#      tracked item, on the Overview sheet (columns E/F) and on each
#      per-language sheet (column C, "Status").
#   2. Shrinks the now-narrower Status column(s) to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-language status columns (zh-cn = E, de-de = F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-language sheets: "Status" column (column C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
